# Automatische test-sync: 2025-08-03 15:15:50
#
# Adds a new "Testmail #20" row to the Logs sheet, a matching
# "Klacht / Probleem" summary row to the Dashboard sheet, updates the
# Dashboard conditional formatting / dimensions and repoints the bar
# chart's category/value series so it covers the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 28
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(28, 1).Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Cells.Item(28, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(28, 3).Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Cells.Item(28, 4).Value = "Klacht / Probleem"
$logs.Cells.Item(28, 5).Value = "Bedankt, we hebben dit doorgestuurd naar klachten@bedrijf.nl."
$logs.Cells.Item(28, 6).Value = "2025-08-03 15:15:14"
$logs.Cells.Item(28, 7).Value = "Ja"
$logs.Cells.Item(28, 8).Value = "Ja"
$logs.Cells.Item(28, 9).Value = "Nee"
$logs.Cells.Item(28, 10).Value = "Nee"

# Extend the existing conditional formatting rules (one per column) so
# they keep covering the full data range, now including row 28.
$logs.Range("D2:D27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D28"))
$logs.Range("G2:G27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G28"))
$logs.Range("H2:H27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H28"))
$logs.Range("I2:I27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I28"))
$logs.Range("J2:J27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J28"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append summary row 7
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(7, 1).Value = "Klacht / Probleem"
$dashboard.Cells.Item(7, 2).Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series to include the new row
# ---------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$7,Dashboard!`$B`$2:`$B`$7,1)"
